$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-09-16 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-09-17 Wednesday", 2) | Out-Null
$d.Content.Find.Execute("89×67=5963", $true, $false, $false, $false, $false, $true, 1, $false, "33×86=2838", 2) | Out-Null
$d.Content.Find.Execute("16×68=1088", $true, $false, $false, $false, $false, $true, 1, $false, "93×20=1860", 2) | Out-Null
$d.Content.Find.Execute("80×27=2160", $true, $false, $false, $false, $false, $true, 1, $false, "35×58=2030", 2) | Out-Null
$d.Content.Find.Execute("24×25=600", $true, $false, $false, $false, $false, $true, 1, $false, "66×20=1320", 2) | Out-Null
$d.Content.Find.Execute("87×62=5394", $true, $false, $false, $false, $false, $true, 1, $false, "60×88=5280", 2) | Out-Null
$d.Content.Find.Execute("29×19=551", $true, $false, $false, $false, $false, $true, 1, $false, "80×45=3600", 2) | Out-Null
$d.Content.Find.Execute("47×67=3149", $true, $false, $false, $false, $false, $true, 1, $false, "87×27=2349", 2) | Out-Null
$d.Content.Find.Execute("20×85=1700", $true, $false, $false, $false, $false, $true, 1, $false, "29×68=1972", 2) | Out-Null
$d.Content.Find.Execute("96×63=6048", $true, $false, $false, $false, $false, $true, 1, $false, "46×98=4508", 2) | Out-Null
$d.Content.Find.Execute("11×24=264", $true, $false, $false, $false, $false, $true, 1, $false, "49×57=2793", 2) | Out-Null
$d.Content.Find.Execute("62×63=3906", $true, $false, $false, $false, $false, $true, 1, $false, "82×94=7708", 2) | Out-Null
$d.Content.Find.Execute("66×93=6138", $true, $false, $false, $false, $false, $true, 1, $false, "44×18=792", 2) | Out-Null
$d.Content.Find.Execute("45×31=1395", $true, $false, $false, $false, $false, $true, 1, $false, "15×28=420", 2) | Out-Null
$d.Content.Find.Execute("77×39=3003", $true, $false, $false, $false, $false, $true, 1, $false, "11×49=539", 2) | Out-Null
$d.Content.Find.Execute("39×36=1404", $true, $false, $false, $false, $false, $true, 1, $false, "48×20=960", 2) | Out-Null
$d.Content.Find.Execute("84×77=6468", $true, $false, $false, $false, $false, $true, 1, $false, "28×82=2296", 2) | Out-Null
$d.Content.Find.Execute("50×93=4650", $true, $false, $false, $false, $false, $true, 1, $false, "79×70=5530", 2) | Out-Null
$d.Content.Find.Execute("83×11=913", $true, $false, $false, $false, $false, $true, 1, $false, "32×45=1440", 2) | Out-Null
$d.Content.Find.Execute("27×34=918", $true, $false, $false, $false, $false, $true, 1, $false, "68×99=6732", 2) | Out-Null
$d.Content.Find.Execute("20×40=800", $true, $false, $false, $false, $false, $true, 1, $false, "82×13=1066", 2) | Out-Null
$d.Content.Find.Execute("23×16=368", $true, $false, $false, $false, $false, $true, 1, $false, "19×16=304", 2) | Out-Null
$d.Content.Find.Execute("15×50=750", $true, $false, $false, $false, $false, $true, 1, $false, "56×75=4200", 2) | Out-Null
$d.Content.Find.Execute("34×34=1156", $true, $false, $false, $false, $false, $true, 1, $false, "51×87=4437", 2) | Out-Null
$d.Content.Find.Execute("78×98=7644", $true, $false, $false, $false, $false, $true, 1, $false, "36×28=1008", 2) | Out-Null
$d.Content.Find.Execute("34×98=3332", $true, $false, $false, $false, $false, $true, 1, $false, "60×61=3660", 2) | Out-Null
